# Update from the master branch
#
# 1) Slide 1, "TextBox 218": second paragraph "UI" -> "레벨"
#    (keep first paragraph "기획서" untouched)
# 2) Slide 3, "사각형: 둥근 모서리 14": "스테이지 1(튜토리얼)" -> "스테이지 1"

$p = $ppt.ActivePresentation

# --- Slide 1: "UI" -> "레벨" -------------------------------------------
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(1)
$tr1 = $sh1.TextFrame.TextRange
$para2 = $tr1.Paragraphs(2, 1)
$para2.Text = "레벨"

# --- Slide 3: "스테이지 1(튜토리얼)" -> "스테이지 1" ---------------------
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(3)
$tr3 = $sh3.TextFrame.TextRange
$tr3.Text = "스테이지 1"
